$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(5920, 45968.95833333334),
    @(5835, 45968.96875),
    @(5744, 45968.97916666666),
    @(5666, 45968.98958333334),
    @(5600, 45969),
    @(5569, 45969.01041666666),
    @(5557, 45969.02083333334),
    @(5475, 45969.03125),
    @(5469, 45969.04166666666),
    @(5449, 45969.05208333334),
    @(5412, 45969.0625),
    @(5365, 45969.07291666666),
    @(5403, 45969.08333333334),
    @(5364, 45969.09375),
    @(5399, 45969.10416666666),
    @(5391, 45969.11458333334),
    @(5423, 45969.125),
    @(5402, 45969.13541666666),
    @(5421, 45969.14583333334),
    @(5401, 45969.15625),
    @(5457, 45969.16666666666),
    @(5475, 45969.17708333334),
    @(5481, 45969.1875),
    @(5526, 45969.19791666666),
    @(5485, 45969.20833333334),
    @(5534, 45969.21875),
    @(5588, 45969.22916666666),
    @(5636, 45969.23958333334),
    @(5744, 45969.25),
    @(5821, 45969.26041666666),
    @(5843, 45969.27083333334),
    @(5915, 45969.28125),
    @(5985, 45969.29166666666),
    @(6024, 45969.30208333334)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $row++
}
